$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2271293375394322
$ws.Range("C2").Value = 0.5078864353312302
$ws.Range("J2").Value = 0.01892744479495268
$ws.Range("P2").Value = 0.167192429022082
$ws.Range("S2").Value = 0.07886435331230283
$ws.Range("B3").Value = 0.006024096385542169
$ws.Range("C3").Value = 0.01807228915662651
$ws.Range("J3").Value = 0.04216867469879518
$ws.Range("P3").Value = 0.6927710843373494
$ws.Range("S3").Value = 0.2409638554216867
$ws.Range("J4").Value = 0.08571428571428572
$ws.Range("P4").Value = 0.6
$ws.Range("S4").Value = 0.3142857142857143
$ws.Range("J5").Value = 0.5
$ws.Range("P5").Value = 0.5
$ws.Range("B6").Value = 0.05429864253393665
$ws.Range("D6").Value = 0.004524886877828055
$ws.Range("F6").Value = 0.08144796380090498
$ws.Range("J6").Value = 0.3031674208144796
$ws.Range("O6").Value = 0.04072398190045249
$ws.Range("Q6").Value = 0.1447963800904978
$ws.Range("R6").Value = 0.04977375565610859
$ws.Range("S6").Value = 0.3212669683257919
$ws.Range("B7").Value = 0.07954545454545454
$ws.Range("D7").Value = 0.01704545454545454
$ws.Range("E7").Value = 0.005681818181818182
$ws.Range("F7").Value = 0.05113636363636364
$ws.Range("J7").Value = 0.1590909090909091
$ws.Range("O7").Value = 0.02840909090909091
$ws.Range("Q7").Value = 0.1534090909090909
$ws.Range("R7").Value = 0.07954545454545454
$ws.Range("S7").Value = 0.4261363636363636
$ws.Range("B8").Value = 0.1044776119402985
$ws.Range("D8").Value = 0.01279317697228145
$ws.Range("F8").Value = 0.05756929637526653
$ws.Range("J8").Value = 0.1130063965884861
$ws.Range("O8").Value = 0.02985074626865672
$ws.Range("Q8").Value = 0.1663113006396588
$ws.Range("R8").Value = 0.06823027718550106
$ws.Range("S8").Value = 0.4477611940298508
$ws.Range("B9").Value = 0.06532663316582915
$ws.Range("D9").Value = 0.01005025125628141
$ws.Range("E9").Value = 0.005025125628140704
$ws.Range("F9").Value = 0.09547738693467336
$ws.Range("J9").Value = 0.1507537688442211
$ws.Range("O9").Value = 0.01005025125628141
$ws.Range("Q9").Value = 0.221105527638191
$ws.Range("R9").Value = 0.07537688442211055
$ws.Range("S9").Value = 0.3668341708542713
$ws.Range("B10").Value = 0.126571668063705
$ws.Range("D10").Value = 0.02095557418273261
$ws.Range("F10").Value = 0.05364626990779547
$ws.Range("J10").Value = 0.1299245599329422
$ws.Range("O10").Value = 0.03017602682313495
$ws.Range("Q10").Value = 0.1860854987426656
$ws.Range("R10").Value = 0.06621961441743504
$ws.Range("S10").Value = 0.3864207879295893
$ws.Range("G11").Value = 0.1835443037974684
$ws.Range("J11").Value = 0.1012658227848101
$ws.Range("K11").Value = 0.2405063291139241
$ws.Range("L11").Value = 0.4367088607594937
$ws.Range("S11").Value = 0.0379746835443038
$ws.Range("G12").Value = 0.7021276595744681
$ws.Range("J12").Value = 0.2127659574468085
$ws.Range("K12").Value = 0.007092198581560284
$ws.Range("L12").Value = 0.01418439716312057
$ws.Range("S12").Value = 0.06382978723404255
$ws.Range("F13").Value = 0.025
$ws.Range("G13").Value = 0.525
$ws.Range("J13").Value = 0.35
$ws.Range("S13").Value = 0.1
$ws.Range("G14").Value = 1
$ws.Range("F15").Value = 0.02631578947368421
$ws.Range("H15").Value = 0.1666666666666667
$ws.Range("I15").Value = 0.09210526315789473
$ws.Range("J15").Value = 0.3026315789473684
$ws.Range("K15").Value = 0.06140350877192982
$ws.Range("M15").Value = 0.0131578947368421
$ws.Range("O15").Value = 0.03947368421052631
$ws.Range("S15").Value = 0.2982456140350877
$ws.Range("F16").Value = 0.0273224043715847
$ws.Range("H16").Value = 0.180327868852459
$ws.Range("I16").Value = 0.09289617486338798
$ws.Range("J16").Value = 0.3551912568306011
$ws.Range("K16").Value = 0.09836065573770492
$ws.Range("M16").Value = 0.0273224043715847
$ws.Range("O16").Value = 0.07103825136612021
$ws.Range("S16").Value = 0.1475409836065574
$ws.Range("F17").Value = 0.03258145363408521
$ws.Range("H17").Value = 0.1904761904761905
$ws.Range("I17").Value = 0.112781954887218
$ws.Range("J17").Value = 0.3784461152882205
$ws.Range("K17").Value = 0.07518796992481203
$ws.Range("M17").Value = 0.007518796992481203
$ws.Range("O17").Value = 0.06015037593984962
$ws.Range("S17").Value = 0.1428571428571428
$ws.Range("F18").Value = 0.02597402597402598
$ws.Range("H18").Value = 0.1233766233766234
$ws.Range("I18").Value = 0.09740259740259741
$ws.Range("J18").Value = 0.435064935064935
$ws.Range("K18").Value = 0.07142857142857142
$ws.Range("M18").Value = 0.03246753246753246
$ws.Range("O18").Value = 0.05844155844155844
$ws.Range("S18").Value = 0.1558441558441558
$ws.Range("F19").Value = 0.02239382239382239
$ws.Range("H19").Value = 0.2378378378378379
$ws.Range("I19").Value = 0.07799227799227799
$ws.Range("J19").Value = 0.3444015444015444
$ws.Range("K19").Value = 0.1196911196911197
$ws.Range("M19").Value = 0.01853281853281853
$ws.Range("N19").Value = 0.001544401544401544
$ws.Range("O19").Value = 0.06254826254826255
$ws.Range("S19").Value = 0.1150579150579151
